$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 93, shifting rows 93:151 down to 94:152.
# Excel duplicates the source row's formatting/content into the new row
# the same way a normal "Insert Copied Cells"/row-insert does, so copy
# row 93 first and insert the copy above itself - this keeps all of the
# unchanged columns (A,B,C,E,F,G,H,I,N,O,Q,R) intact automatically.
$ws.Rows.Item(93).Copy()
$ws.Rows.Item(93).Insert()

# Now overwrite the columns that actually differ for the new record
$ws.Cells.Item(93, 4).Value = 44596   # D93 date
$ws.Cells.Item(93, 10).Value = 34     # J93 volumen
$ws.Cells.Item(93, 11).Value = 8000   # K93 precio minimo
$ws.Cells.Item(93, 12).Value = 9000   # L93 precio maximo
$ws.Cells.Item(93, 13).Value = 8441   # M93 precio promedio ponderado
$ws.Cells.Item(93, 16).Value = 2814   # P93 precio $/Kg
